# Update NATMI LR-pair results (Il1rn-Il1r1) with recomputed TPM-based values.
# Columns A-D (cluster/gene labels) and K-L (receptor-expressing cells / detection rate)
# are unchanged; columns E-J (ligand side) and M-T (receptor side + edge weights) are updated
# for data rows 2-25, matching the new TPM input used to regenerate the NATMI output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = New-Object 'double[,]' 24,16
$newValues[0,0] = 3
$newValues[0,1] = 1
$newValues[0,2] = [double]"1.498484"
$newValues[0,3] = [double]"4.495452"
$newValues[0,4] = [double]"0.001216692661445581"
$newValues[0,5] = [double]"0.001216692661445581"
$newValues[0,6] = 3
$newValues[0,7] = 1
$newValues[0,8] = [double]"16.087096"
$newValues[0,9] = [double]"48.261288"
$newValues[0,10] = [double]"0.1263055268415452"
$newValues[0,11] = [double]"0.1263055268415452"
$newValues[0,12] = [double]"24.106255962464"
$newValues[0,13] = [double]"216.956303662176"
$newValues[0,14] = [double]"0.0001536750076081259"
$newValues[0,15] = [double]"0.0001536750076081259"
$newValues[1,0] = 3
$newValues[1,1] = 1
$newValues[1,2] = [double]"1.498484"
$newValues[1,3] = [double]"4.495452"
$newValues[1,4] = [double]"0.001216692661445581"
$newValues[1,5] = [double]"0.001216692661445581"
$newValues[1,6] = 3
$newValues[1,7] = 1
$newValues[1,8] = [double]"95.39815"
$newValues[1,9] = [double]"286.19445"
$newValues[1,10] = [double]"0.7490048915888087"
$newValues[1,11] = [double]"0.7490048915888088"
$newValues[1,12] = [double]"142.9526014046"
$newValues[1,13] = [double]"1286.5734126414"
$newValues[1,14] = [double]"0.0009113087549829462"
$newValues[1,15] = [double]"0.0009113087549829463"
$newValues[2,0] = 3
$newValues[2,1] = 1
$newValues[2,2] = [double]"1.498484"
$newValues[2,3] = [double]"4.495452"
$newValues[2,4] = [double]"0.001216692661445581"
$newValues[2,5] = [double]"0.001216692661445581"
$newValues[2,6] = 3
$newValues[2,7] = 1
$newValues[2,8] = [double]"0.5200936666666667"
$newValues[2,9] = [double]"1.560281"
$newValues[2,10] = [double]"0.004083440825819921"
$newValues[2,11] = [double]"0.004083440825819921"
$newValues[2,12] = [double]"0.7793520380013335"
$newValues[2,13] = [double]"7.014168342012001"
$newValues[2,14] = [double]"4.968292486222379e-06"
$newValues[2,15] = [double]"4.968292486222379e-06"
$newValues[3,0] = 3
$newValues[3,1] = 1
$newValues[3,2] = [double]"1.498484"
$newValues[3,3] = [double]"4.495452"
$newValues[3,4] = [double]"0.001216692661445581"
$newValues[3,5] = [double]"0.001216692661445581"
$newValues[3,6] = 3
$newValues[3,7] = 1
$newValues[3,8] = [double]"14.15205133333333"
$newValues[3,9] = [double]"42.456154"
$newValues[3,10] = [double]"0.1111128011883101"
$newValues[3,11] = [double]"0.1111128011883101"
$newValues[3,12] = [double]"21.20662249017867"
$newValues[3,13] = [double]"190.859602411608"
$newValues[3,14] = [double]"0.0001351901297984787"
$newValues[3,15] = [double]"0.0001351901297984787"
$newValues[4,0] = 3
$newValues[4,1] = 1
$newValues[4,2] = [double]"1.498484"
$newValues[4,3] = [double]"4.495452"
$newValues[4,4] = [double]"0.001216692661445581"
$newValues[4,5] = [double]"0.001216692661445581"
$newValues[4,6] = 3
$newValues[4,7] = 1
$newValues[4,8] = [double]"0.794831"
$newValues[4,9] = [double]"2.384493"
$newValues[4,10] = [double]"0.006240501592393819"
$newValues[4,11] = [double]"0.006240501592393819"
$newValues[4,12] = [double]"1.191041536204"
$newValues[4,13] = [double]"10.719373825836"
$newValues[4,14] = [double]"7.592772491205018e-06"
$newValues[4,15] = [double]"7.592772491205018e-06"
$newValues[5,0] = 3
$newValues[5,1] = 1
$newValues[5,2] = [double]"1.498484"
$newValues[5,3] = [double]"4.495452"
$newValues[5,4] = [double]"0.001216692661445581"
$newValues[5,5] = [double]"0.001216692661445581"
$newValues[5,6] = 3
$newValues[5,7] = 1
$newValues[5,8] = [double]"0.4143026666666667"
$newValues[5,9] = [double]"1.242908"
$newValues[5,10] = [double]"0.003252837963122146"
$newValues[5,11] = [double]"0.003252837963122146"
$newValues[5,12] = [double]"0.6208259171573334"
$newValues[5,13] = [double]"5.587433254416"
$newValues[5,14] = [double]"3.957704078602306e-06"
$newValues[5,15] = [double]"3.957704078602306e-06"
$newValues[6,0] = 3
$newValues[6,1] = 1
$newValues[6,2] = [double]"140.872935"
$newValues[6,3] = [double]"422.618805"
$newValues[6,4] = [double]"0.1143816458572799"
$newValues[6,5] = [double]"0.1143816458572799"
$newValues[6,6] = 3
$newValues[6,7] = 1
$newValues[6,8] = [double]"16.087096"
$newValues[6,9] = [double]"48.261288"
$newValues[6,10] = [double]"0.1263055268415452"
$newValues[6,11] = [double]"0.1263055268415452"
$newValues[6,12] = [double]"2266.23642914676"
$newValues[6,13] = [double]"20396.12786232084"
$newValues[6,14] = [double]"0.01444703404100679"
$newValues[6,15] = [double]"0.01444703404100679"
$newValues[7,0] = 3
$newValues[7,1] = 1
$newValues[7,2] = [double]"140.872935"
$newValues[7,3] = [double]"422.618805"
$newValues[7,4] = [double]"0.1143816458572799"
$newValues[7,5] = [double]"0.1143816458572799"
$newValues[7,6] = 3
$newValues[7,7] = 1
$newValues[7,8] = [double]"95.39815"
$newValues[7,9] = [double]"286.19445"
$newValues[7,10] = [double]"0.7490048915888087"
$newValues[7,11] = [double]"0.7490048915888088"
$newValues[7,12] = [double]"13439.01738407025"
$newValues[7,13] = [double]"120951.1564566323"
$newValues[7,14] = [double]"0.08567241225508146"
$newValues[7,15] = [double]"0.08567241225508146"
$newValues[8,0] = 3
$newValues[8,1] = 1
$newValues[8,2] = [double]"140.872935"
$newValues[8,3] = [double]"422.618805"
$newValues[8,4] = [double]"0.1143816458572799"
$newValues[8,5] = [double]"0.1143816458572799"
$newValues[8,6] = 3
$newValues[8,7] = 1
$newValues[8,8] = [double]"0.5200936666666667"
$newValues[8,9] = [double]"1.560281"
$newValues[8,10] = [double]"0.004083440825819921"
$newValues[8,11] = [double]"0.004083440825819921"
$newValues[8,12] = [double]"73.267121298245"
$newValues[8,13] = [double]"659.404091684205"
$newValues[8,14] = [double]"0.0004670706824180929"
$newValues[8,15] = [double]"0.0004670706824180928"
$newValues[9,0] = 3
$newValues[9,1] = 1
$newValues[9,2] = [double]"140.872935"
$newValues[9,3] = [double]"422.618805"
$newValues[9,4] = [double]"0.1143816458572799"
$newValues[9,5] = [double]"0.1143816458572799"
$newValues[9,6] = 3
$newValues[9,7] = 1
$newValues[9,8] = [double]"14.15205133333333"
$newValues[9,9] = [double]"42.456154"
$newValues[9,10] = [double]"0.1111128011883101"
$newValues[9,11] = [double]"0.1111128011883101"
$newValues[9,12] = [double]"1993.64100759733"
$newValues[9,13] = [double]"17942.76906837597"
$newValues[9,14] = [double]"0.01270926507573164"
$newValues[9,15] = [double]"0.01270926507573164"
$newValues[10,0] = 3
$newValues[10,1] = 1
$newValues[10,2] = [double]"140.872935"
$newValues[10,3] = [double]"422.618805"
$newValues[10,4] = [double]"0.1143816458572799"
$newValues[10,5] = [double]"0.1143816458572799"
$newValues[10,6] = 3
$newValues[10,7] = 1
$newValues[10,8] = [double]"0.794831"
$newValues[10,9] = [double]"2.384493"
$newValues[10,10] = [double]"0.006240501592393819"
$newValues[10,11] = [double]"0.006240501592393819"
$newValues[10,12] = [double]"111.970175798985"
$newValues[10,13] = [double]"1007.731582190865"
$newValues[10,14] = [double]"0.0007137988431129812"
$newValues[10,15] = [double]"0.0007137988431129812"
$newValues[11,0] = 3
$newValues[11,1] = 1
$newValues[11,2] = [double]"140.872935"
$newValues[11,3] = [double]"422.618805"
$newValues[11,4] = [double]"0.1143816458572799"
$newValues[11,5] = [double]"0.1143816458572799"
$newValues[11,6] = 3
$newValues[11,7] = 1
$newValues[11,8] = [double]"0.4143026666666667"
$newValues[11,9] = [double]"1.242908"
$newValues[11,10] = [double]"0.003252837963122146"
$newValues[11,11] = [double]"0.003252837963122146"
$newValues[11,12] = [double]"58.36403263166"
$newValues[11,13] = [double]"525.27629368494"
$newValues[11,14] = [double]"0.0003720649599289531"
$newValues[11,15] = [double]"0.0003720649599289531"
$newValues[12,0] = 3
$newValues[12,1] = 1
$newValues[12,2] = [double]"1017.335652666667"
$newValues[12,3] = [double]"3052.006958"
$newValues[12,4] = [double]"0.8260247175321745"
$newValues[12,5] = [double]"0.8260247175321747"
$newValues[12,6] = 3
$newValues[12,7] = 1
$newValues[12,8] = [double]"16.087096"
$newValues[12,9] = [double]"48.261288"
$newValues[12,10] = [double]"0.1263055268415452"
$newValues[12,11] = [double]"0.1263055268415452"
$newValues[12,12] = [double]"16365.97630867132"
$newValues[12,13] = [double]"147293.7867780419"
$newValues[12,14] = [double]"0.1043314871320399"
$newValues[12,15] = [double]"0.1043314871320399"
$newValues[13,0] = 3
$newValues[13,1] = 1
$newValues[13,2] = [double]"1017.335652666667"
$newValues[13,3] = [double]"3052.006958"
$newValues[13,4] = [double]"0.8260247175321745"
$newValues[13,5] = [double]"0.8260247175321747"
$newValues[13,6] = 3
$newValues[13,7] = 1
$newValues[13,8] = [double]"95.39815"
$newValues[13,9] = [double]"286.19445"
$newValues[13,10] = [double]"0.7490048915888087"
$newValues[13,11] = [double]"0.7490048915888088"
$newValues[13,12] = [double]"97051.93919344257"
$newValues[13,13] = [double]"873467.4527409831"
$newValues[13,14] = [double]"0.6186965540048627"
$newValues[13,15] = [double]"0.6186965540048629"
$newValues[14,0] = 3
$newValues[14,1] = 1
$newValues[14,2] = [double]"1017.335652666667"
$newValues[14,3] = [double]"3052.006958"
$newValues[14,4] = [double]"0.8260247175321745"
$newValues[14,5] = [double]"0.8260247175321747"
$newValues[14,6] = 3
$newValues[14,7] = 1
$newValues[14,8] = [double]"0.5200936666666667"
$newValues[14,9] = [double]"1.560281"
$newValues[14,10] = [double]"0.004083440825819921"
$newValues[14,11] = [double]"0.004083440825819921"
$newValues[14,12] = [double]"529.1098298261331"
$newValues[14,13] = [double]"4761.988468435198"
$newValues[14,14] = [double]"0.00337302305470725"
$newValues[14,15] = [double]"0.00337302305470725"
$newValues[15,0] = 3
$newValues[15,1] = 1
$newValues[15,2] = [double]"1017.335652666667"
$newValues[15,3] = [double]"3052.006958"
$newValues[15,4] = [double]"0.8260247175321745"
$newValues[15,5] = [double]"0.8260247175321747"
$newValues[15,6] = 3
$newValues[15,7] = 1
$newValues[15,8] = [double]"14.15205133333333"
$newValues[15,9] = [double]"42.456154"
$newValues[15,10] = [double]"0.1111128011883101"
$newValues[15,11] = [double]"0.1111128011883101"
$newValues[15,12] = [double]"14397.38637976884"
$newValues[15,13] = [double]"129576.4774179195"
$newValues[15,14] = [double]"0.09178192021578255"
$newValues[15,15] = [double]"0.09178192021578256"
$newValues[16,0] = 3
$newValues[16,1] = 1
$newValues[16,2] = [double]"1017.335652666667"
$newValues[16,3] = [double]"3052.006958"
$newValues[16,4] = [double]"0.8260247175321745"
$newValues[16,5] = [double]"0.8260247175321747"
$newValues[16,6] = 3
$newValues[16,7] = 1
$newValues[16,8] = [double]"0.794831"
$newValues[16,9] = [double]"2.384493"
$newValues[16,10] = [double]"0.006240501592393819"
$newValues[16,11] = [double]"0.006240501592393819"
$newValues[16,12] = [double]"808.6099141446992"
$newValues[16,13] = [double]"7277.489227302293"
$newValues[16,14] = [double]"0.005154808565116189"
$newValues[16,15] = [double]"0.00515480856511619"
$newValues[17,0] = 3
$newValues[17,1] = 1
$newValues[17,2] = [double]"1017.335652666667"
$newValues[17,3] = [double]"3052.006958"
$newValues[17,4] = [double]"0.8260247175321745"
$newValues[17,5] = [double]"0.8260247175321747"
$newValues[17,6] = 3
$newValues[17,7] = 1
$newValues[17,8] = [double]"0.4143026666666667"
$newValues[17,9] = [double]"1.242908"
$newValues[17,10] = [double]"0.003252837963122146"
$newValues[17,11] = [double]"0.003252837963122146"
$newValues[17,12] = [double]"421.4848737948737"
$newValues[17,13] = [double]"3793.363864153864"
$newValues[17,14] = [double]"0.002686924559665905"
$newValues[17,15] = [double]"0.002686924559665905"
$newValues[18,0] = 3
$newValues[18,1] = 1
$newValues[18,2] = [double]"71.89729933333332"
$newValues[18,3] = [double]"215.691898"
$newValues[18,4] = [double]"0.0583769439490998"
$newValues[18,5] = [double]"0.05837694394909981"
$newValues[18,6] = 3
$newValues[18,7] = 1
$newValues[18,8] = [double]"16.087096"
$newValues[18,9] = [double]"48.261288"
$newValues[18,10] = [double]"0.1263055268415452"
$newValues[18,11] = [double]"0.1263055268415452"
$newValues[18,12] = [double]"1156.618756516069"
$newValues[18,13] = [double]"10409.56880864462"
$newValues[18,14] = [double]"0.007373330660890406"
$newValues[18,15] = [double]"0.007373330660890407"
$newValues[19,0] = 3
$newValues[19,1] = 1
$newValues[19,2] = [double]"71.89729933333332"
$newValues[19,3] = [double]"215.691898"
$newValues[19,4] = [double]"0.0583769439490998"
$newValues[19,5] = [double]"0.05837694394909981"
$newValues[19,6] = 3
$newValues[19,7] = 1
$newValues[19,8] = [double]"95.39815"
$newValues[19,9] = [double]"286.19445"
$newValues[19,10] = [double]"0.7490048915888087"
$newValues[19,11] = [double]"0.7490048915888088"
$newValues[19,12] = [double]"6858.869346396233"
$newValues[19,13] = [double]"61729.8241175661"
$newValues[19,14] = [double]"0.04372461657388146"
$newValues[19,15] = [double]"0.04372461657388146"
$newValues[20,0] = 3
$newValues[20,1] = 1
$newValues[20,2] = [double]"71.89729933333332"
$newValues[20,3] = [double]"215.691898"
$newValues[20,4] = [double]"0.0583769439490998"
$newValues[20,5] = [double]"0.05837694394909981"
$newValues[20,6] = 3
$newValues[20,7] = 1
$newValues[20,8] = [double]"0.5200936666666667"
$newValues[20,9] = [double]"1.560281"
$newValues[20,10] = [double]"0.004083440825819921"
$newValues[20,11] = [double]"0.004083440825819921"
$newValues[20,12] = [double]"37.39333003370422"
$newValues[20,13] = [double]"336.539970303338"
$newValues[20,14] = [double]"0.0002383787962083553"
$newValues[20,15] = [double]"0.0002383787962083554"
$newValues[21,0] = 3
$newValues[21,1] = 1
$newValues[21,2] = [double]"71.89729933333332"
$newValues[21,3] = [double]"215.691898"
$newValues[21,4] = [double]"0.0583769439490998"
$newValues[21,5] = [double]"0.05837694394909981"
$newValues[21,6] = 3
$newValues[21,7] = 1
$newValues[21,8] = [double]"14.15205133333333"
$newValues[21,9] = [double]"42.456154"
$newValues[21,10] = [double]"0.1111128011883101"
$newValues[21,11] = [double]"0.1111128011883101"
$newValues[21,12] = [double]"1017.494270893366"
$newValues[21,13] = [double]"9157.44843804029"
$newValues[21,14] = [double]"0.00648642576699745"
$newValues[21,15] = [double]"0.00648642576699745"
$newValues[22,0] = 3
$newValues[22,1] = 1
$newValues[22,2] = [double]"71.89729933333332"
$newValues[22,3] = [double]"215.691898"
$newValues[22,4] = [double]"0.0583769439490998"
$newValues[22,5] = [double]"0.05837694394909981"
$newValues[22,6] = 3
$newValues[22,7] = 1
$newValues[22,8] = [double]"0.794831"
$newValues[22,9] = [double]"2.384493"
$newValues[22,10] = [double]"0.006240501592393819"
$newValues[22,11] = [double]"0.006240501592393819"
$newValues[22,12] = [double]"57.14620232641266"
$newValues[22,13] = [double]"514.315820937714"
$newValues[22,14] = [double]"0.000364301411673442"
$newValues[22,15] = [double]"0.0003643014116734421"
$newValues[23,0] = 3
$newValues[23,1] = 1
$newValues[23,2] = [double]"71.89729933333332"
$newValues[23,3] = [double]"215.691898"
$newValues[23,4] = [double]"0.0583769439490998"
$newValues[23,5] = [double]"0.05837694394909981"
$newValues[23,6] = 3
$newValues[23,7] = 1
$newValues[23,8] = [double]"0.4143026666666667"
$newValues[23,9] = [double]"1.242908"
$newValues[23,10] = [double]"0.003252837963122146"
$newValues[23,11] = [double]"0.003252837963122146"
$newValues[23,12] = [double]"29.78724283993155"
$newValues[23,13] = [double]"268.0851855593839"
$newValues[23,14] = [double]"0.0001898907394486855"
$newValues[23,15] = [double]"0.0001898907394486855"

$ws.Range("E2:T25").Value2 = $newValues

